$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21; existing rows 21-95 shift down to 22-96.
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly data point.
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44414
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100112037
$ws.Range("G21").Value = "Cebollín"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 3200
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 1000
$ws.Range("M21").Value = 950
$ws.Range("N21").Value = "$/paquete 6 unidades"
$ws.Range("O21").Value = "Provincia del Elquí"
$ws.Range("P21").Value = 158
$ws.Range("Q21").Value = 6
$ws.Range("R21").Value = "Hortaliza"
